$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Description column (H) for the LED rows to append the
# "0201" package designation (and drop the "(0603 Metric)" suffix on the
# Blue LED row).
$ws.Range("H36").Value = "Red 631nm LED Indication - Discrete 2.4V 2-SMD, No Lead 0201"
$ws.Range("H37").Value = "Yellow 589nm LED Indication - Discrete 2.4V 2-SMD, No Lead 0201"
$ws.Range("H38").Value = "Green 571nm LED Indication - Discrete 2.4V 2-SMD, No Lead 0201"
$ws.Range("H39").Value = "Blue 468nm LED Indication - Discrete 2.9V 0201 "

# Restore the view state: scrolled so column C is the left-most visible
# column, with H29 as the active selected cell.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("H29").Select()
